# Insert a new data row at position 274 (pushing existing rows 274:387 down to
# 275:388) and populate it with the new Betarraga price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("274:274").Insert()

$ws.Range("A274").Value = 10
$ws.Range("B274").Value = "Vega Modelo de Temuco"
$ws.Range("C274").Value = "La Araucanía"
$ws.Range("D274").Value = 44704
$ws.Range("E274").Value = 9
$ws.Range("F274").Value = 100114014
$ws.Range("G274").Value = "Betarraga"
$ws.Range("H274").Value = "Sin especificar"
$ws.Range("I274").Value = "Primera"
$ws.Range("J274").Value = 40
$ws.Range("K274").Value = 8000
$ws.Range("L274").Value = 8000
$ws.Range("M274").Value = 8000
$ws.Range("N274").Value = '$/docena de paquetes'
$ws.Range("O274").Value = "Provincia de Cautín"
$ws.Range("P274").Value = 667
$ws.Range("Q274").Value = 12
$ws.Range("R274").Value = "Hortaliza"
